$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("dob") on rows 2 and 4 currently hold a numeric Excel date
# serial with a date number-format applied. Replace with the literal text
# value "05/10/2022" (leading apostrophe forces text, like typing it into
# Excel) and strip the date formatting back off the cell.
$ws.Range("G2").Value = "'05/10/2022"
$ws.Range("G2").ClearFormats()

$ws.Range("G4").Value = "'05/10/2022"
$ws.Range("G4").ClearFormats()

# Columns D ("email") / E ("contactNumber") become empty text (not blank
# cells) - a lone leading apostrophe enters an empty text value, same as
# typing "'" then nothing into Excel.
$ws.Range("E2").Value = "'"
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "'"
$ws.Range("E3").ClearFormats()

$ws.Range("D4").Value = "'"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = "'"
$ws.Range("E4").ClearFormats()

# Column I ("staff"): row 2 and row 3's "[]" become an empty string, and
# row 4's "[1 2]" becomes "1;2".
$ws.Range("I2").Value = "'"
$ws.Range("I2").ClearFormats()

$ws.Range("I3").Value = "'"
$ws.Range("I3").ClearFormats()

$ws.Range("I4").Value = "1;2"
